# Updates cryptos list data (prices and 1h volume %) as scraped by GitHub Actions.
# Also corrects two mis-ordered rows (WrappedBTC/Polkadot and Dai/RenderToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price strings are written with a leading apostrophe and then
# restored to the default style so Excel keeps them as text (preserving
# formatting such as trailing zeros) instead of silently converting them to
# floating point numbers.

$ws.Range("D2").Value = '62.224.40'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '2.999.53'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'580.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").Value = "'146.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.09%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.002.23'
$ws.Range("E8").Value = '  -2.09%  '
$ws.Range("D9").Value = "'0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.20%  '
$ws.Range("D10").Value = "'0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.78%  '
$ws.Range("D11").Value = "'5.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.14%  '
$ws.Range("D12").Value = "'0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = "'0.0000228"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.72%  '
$ws.Range("D14").Value = "'34.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.08%  '
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").Value = '3.506.91'
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '62.291.79'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'7.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").Value = '3.013.77'
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").Value = "'459.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.09%  '
$ws.Range("D21").Value = "'13.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").Value = "'0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.85%  '
$ws.Range("D23").Value = "'7.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("D24").Value = "'2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.91%  '
$ws.Range("D25").Value = "'79.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").Value = "'12.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.70%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = "'9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.62%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").Value = "'2.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("D31").Value = "'7.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.74%  '
$ws.Range("D32").Value = "'2.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.38%  '
$ws.Range("D33").Value = "'27.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = "'0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.35%  '
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("D36").Value = '0.0₃0779'
$ws.Range("E36").Value = '  -4.63%  '
$ws.Range("D37").Value = "'5.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.16%  '
$ws.Range("D38").Value = "'2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.01%  '
$ws.Range("D39").Value = "'50.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").Value = "'8.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("D41").Value = "'2.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.72%  '
$ws.Range("D42").Value = "'414.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.99%  '
$ws.Range("D43").Value = "'0.112"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("D44").Value = "'0.273"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.17%  '
$ws.Range("D45").Value = '2.761.54'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("D46").Value = "'0.0350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("D47").Value = "'38.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.57%  '
$ws.Range("D48").Value = "'128.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D50").Value = "'0.107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.10%  '
$ws.Range("D51").Value = "'23.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.67%  '
